# Update the timestamp in column A for rows 2-13 on the "ランサーズ" sheet
# from "2025-09-08 12:36:36" to "2025-09-08 12:48:46".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-08 12:48:46"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
